# AFDP-157 - Add Access Control List to SOLR documents - apply assignment
# and data access control rules to case files.
#
# The "Set Assignee" rule row (row 21 of the Save Case Rules table, which
# hard-coded every case to be assigned to 'ann-acm') is removed from the
# rule sheet now that assignment is handled by the new access-control
# rules instead of this rule table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire "Set Assignee" row - everything below shifts up.
$ws.Rows.Item(21).Delete()

# Leave the selection where the edit happened.
$ws.Range("C22").Select() | Out-Null

# Match the saved window chrome state from the commit (cosmetic; some
# hosts may not persist this view-only attribute).
$excel.ActiveWindow().TabRatio = 0.134
